$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.112.43"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "3.424.28"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'409.62"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'129.47"
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +6.62%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.746"
$ws.Range("E9").Value = "  +10.57%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  +19.91%  "
$ws.Range("D11").Value = "'42.71"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("E12").Value = "  +72.47%  "
$ws.Range("D14").Value = "3.964.23"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "'8.97"
$ws.Range("E15").Value = "  +6.16%  "
$ws.Range("D16").Value = "'21.06"
$ws.Range("E16").Value = "  +5.77%  "
$ws.Range("D17").Value = "3.419.45"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "'12.41"
$ws.Range("E18").Value = "  +12.91%  "
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("D20").Value = "62.065.82"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'403.67"
$ws.Range("E21").Value = "  +27.19%  "
$ws.Range("D22").Value = "'89.71"
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'13.24"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").Value = "'32.84"
$ws.Range("E26").Value = "  +10.23%  "
$ws.Range("D27").Value = "'8.75"
$ws.Range("E27").Value = "  +6.05%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.63"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.71"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'11.89"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").Value = "'43.11"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("D37").Value = "'53.91"
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("E40").Value = "  +7.06%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("E42").Value = "  +7.02%  "
$ws.Range("D43").Value = "'141.87"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("D45").Value = "'4.11"
$ws.Range("E45").Value = "  +3.30%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  +8.82%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "'21.90"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("D49").Value = "2.122.57"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'2.38"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.01"
$ws.Range("E51").Value = "  +7.05%  "
